$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4002801.2
$ws.Range("I132").Value = 5266056
$ws.Range("J132").Value = 2494.9167
$ws.Range("K132").Value = 15798168
$ws.Range("L132").Value = 7484.750100000001
$ws.Range("M132").Value = -15795638
$ws.Range("N132").Value = -12544.7501
$ws.Range("H137").Value = 3172.14
$ws.Range("I137").Value = 3049.0444
$ws.Range("J137").Value = 4280
$ws.Range("K137").Value = 9147.1332
$ws.Range("L137").Value = 12840
$ws.Range("M137").Value = -6597.1332
$ws.Range("N137").Value = -17940

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 40000
$ws.Range("J24").Value = 40000
$ws.Range("L24").Value = 40000
$ws.Range("N24").Value = -40748
$ws.Range("H32").Value = 1057.1
$ws.Range("I32").Value = 1063.0101
$ws.Range("J32").Value = 472
$ws.Range("K32").Value = 1063.0101
$ws.Range("L32").Value = 472
$ws.Range("M32").Value = -776.0101
$ws.Range("N32").Value = -1046
$ws.Range("H45").Value = 1399.8379
$ws.Range("I45").Value = 1073.4814
$ws.Range("J45").Value = 2281
$ws.Range("K45").Value = 1073.4814
$ws.Range("L45").Value = 2281
$ws.Range("M45").Value = -696.4813999999999
$ws.Range("N45").Value = -3035
$ws.Range("H61").Value = 2330.3684
$ws.Range("I61").Value = 1086.4
$ws.Range("J61").Value = 3712.5557
$ws.Range("K61").Value = 1086.4
$ws.Range("L61").Value = 3712.5557
$ws.Range("M61").Value = -874.4000000000001
$ws.Range("N61").Value = -4136.5557
$ws.Range("H95").Value = 29808
$ws.Range("J95").Value = 29808
$ws.Range("L95").Value = 29808
$ws.Range("N95").Value = -35300
$ws.Range("H100").Value = 40000
$ws.Range("J100").Value = 40000
$ws.Range("L100").Value = 40000
$ws.Range("N100").Value = -42164
$ws.Range("H122").Value = 3267.7144
$ws.Range("I122").Value = 2566
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 7698
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -5248
$ws.Range("N122").Value = -23650
$ws.Range("H132").Value = 2078.738
$ws.Range("I132").Value = 1556.4482
$ws.Range("J132").Value = 3243.8462
$ws.Range("K132").Value = 4669.3446
$ws.Range("L132").Value = 9731.5386
$ws.Range("M132").Value = -2139.3446
$ws.Range("N132").Value = -14791.5386
$ws.Range("H136").Value = 2330.3684
$ws.Range("I136").Value = 1086.4
$ws.Range("J136").Value = 3712.5557
$ws.Range("K136").Value = 3259.2
$ws.Range("L136").Value = 11137.6671
$ws.Range("M136").Value = -709.2000000000003
$ws.Range("N136").Value = -16237.6671
$ws.Range("H139").Value = 28940
$ws.Range("J139").Value = 28940
$ws.Range("L139").Value = 28940
$ws.Range("N139").Value = -39220

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 10096.5
$ws.Range("J81").Value = 10096.5
$ws.Range("L81").Value = 10096.5
$ws.Range("N81").Value = -12218.5
$ws.Range("H84").Value = 10096.5
$ws.Range("J84").Value = 10096.5
$ws.Range("L84").Value = 30289.5
$ws.Range("N84").Value = -40897.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2909.66
$ws.Range("I31").Value = 1774.5151
$ws.Range("J31").Value = 5113.1763
$ws.Range("K31").Value = 1774.5151
$ws.Range("L31").Value = 5113.1763
$ws.Range("M31").Value = -1479.5151
$ws.Range("N31").Value = -5703.1763
$ws.Range("H34").Value = 2909.66
$ws.Range("I34").Value = 1774.5151
$ws.Range("J34").Value = 5113.1763
$ws.Range("K34").Value = 1774.5151
$ws.Range("L34").Value = 5113.1763
$ws.Range("M34").Value = -1572.5151
$ws.Range("N34").Value = -5517.1763
$ws.Range("H99").Value = 2942.7144
$ws.Range("I99").Value = 1666.6666
$ws.Range("J99").Value = 3899.75
$ws.Range("K99").Value = 1666.6666
$ws.Range("L99").Value = 3899.75
$ws.Range("M99").Value = -168.6666
$ws.Range("N99").Value = -6895.75
$ws.Range("H126").Value = 2942.7144
$ws.Range("I126").Value = 1666.6666
$ws.Range("J126").Value = 3899.75
$ws.Range("K126").Value = 4999.9998
$ws.Range("L126").Value = 11699.25
$ws.Range("M126").Value = -2529.9998
$ws.Range("N126").Value = -16639.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 10558
$ws.Range("I11").Value = 270
$ws.Range("J11").Value = 17416.666
$ws.Range("K11").Value = 810
$ws.Range("L11").Value = 52249.99800000001
$ws.Range("M11").Value = -670
$ws.Range("N11").Value = -52529.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4328.5186
$ws.Range("I70").Value = 4465.55
$ws.Range("J70").Value = 3937
$ws.Range("K70").Value = 4465.55
$ws.Range("L70").Value = 3937
$ws.Range("M70").Value = -4195.55
$ws.Range("N70").Value = -4477
$ws.Range("H73").Value = 4328.5186
$ws.Range("I73").Value = 4465.55
$ws.Range("J73").Value = 3937
$ws.Range("K73").Value = 4465.55
$ws.Range("L73").Value = 3937
$ws.Range("M73").Value = -3529.55
$ws.Range("N73").Value = -5809
$ws.Range("H102").Value = 26135.715
$ws.Range("I102").Value = 1653.7407
$ws.Range("J102").Value = 70203.266
$ws.Range("K102").Value = 1653.7407
$ws.Range("L102").Value = 70203.266
$ws.Range("M102").Value = -31.74070000000006
$ws.Range("N102").Value = -73447.266
$ws.Range("H132").Value = 3114.3777
$ws.Range("I132").Value = 2907.0571
$ws.Range("J132").Value = 3840
$ws.Range("K132").Value = 8721.1713
$ws.Range("L132").Value = 11520
$ws.Range("M132").Value = -6191.1713
$ws.Range("N132").Value = -16580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 238584.81
$ws.Range("I14").Value = 836801
$ws.Range("J14").Value = 14253.75
$ws.Range("K14").Value = 836801
$ws.Range("L14").Value = 14253.75
$ws.Range("M14").Value = -836629
$ws.Range("N14").Value = -14597.75
$ws.Range("H43").Value = 50009.332
$ws.Range("J43").Value = 50009.332
$ws.Range("L43").Value = 50009.332
$ws.Range("N43").Value = -50395.332
$ws.Range("H61").Value = 200005600
$ws.Range("I61").Value = 500001000
$ws.Range("J61").Value = 8668.333000000001
$ws.Range("K61").Value = 500001000
$ws.Range("L61").Value = 8668.333000000001
$ws.Range("M61").Value = -500000798
$ws.Range("N61").Value = -9072.333000000001
$ws.Range("H113").Value = 200005600
$ws.Range("I113").Value = 500001000
$ws.Range("J113").Value = 8668.333000000001
$ws.Range("K113").Value = 500001000
$ws.Range("L113").Value = 8668.333000000001
$ws.Range("M113").Value = -499998830
$ws.Range("N113").Value = -13008.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 661.7857
$ws.Range("I107").Value = 297.14285
$ws.Range("K107").Value = 891.4285500000001
$ws.Range("M107").Value = 1028.57145
$ws.Range("H136").Value = 913.31915
$ws.Range("I136").Value = 479.96667
$ws.Range("J136").Value = 1678.0588
$ws.Range("K136").Value = 1479.90001
$ws.Range("L136").Value = 5034.1764
$ws.Range("M136").Value = 1110.09999
$ws.Range("N136").Value = -10134.1764
